$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 30
$ws1.Range("F3").Value = 178
$ws1.Range("F7").Value = 1723
$ws1.Range("F8").Value = 37
$ws1.Range("F11").Value = 1813
$ws1.Range("F14").Value = 426
$ws1.Range("F22").Value = 771

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 30
$ws4.Range("F3").Value = 178
$ws4.Range("F7").Value = 1723
$ws4.Range("F9").Value = 37
$ws4.Range("F12").Value = 1813
$ws4.Range("F15").Value = 426
$ws4.Range("F23").Value = 771
